# Initial Data File Update
# - Fix category of row 175 (Golosina -> Golosinas)
# - Append 4 new transaction rows (178-181) to "Transacciones"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# --- Fix D175 category -------------------------------------------------
$ws.Range("D175").Value = "Golosinas"

# --- Helper: clone the formatting of the last existing row (177) down to
#     a new row, without touching the values, so the new rows inherit the
#     same look (date style, "disponible" style, etc.) as the rest of the
#     table. -----------------------------------------------------------
function Copy-RowFormat($fromRow, $toRow) {
    $ws.Range("A$fromRow").Copy() | Out-Null
    $ws.Range("A$toRow").PasteSpecial(-4122) | Out-Null
    $ws.Range("P$fromRow").Copy() | Out-Null
    $ws.Range("P$toRow").PasteSpecial(-4122) | Out-Null
}

Copy-RowFormat 177 178
Copy-RowFormat 177 179
Copy-RowFormat 177 180
Copy-RowFormat 177 181

# --- Row 178 -------------------------------------------------------------
$ws.Range("A178").Value = 43601
$ws.Range("B178").Value = 19.5
$ws.Range("C178").Value = "Panquecitos"
$ws.Range("D178").Value = "Golosinas"
$ws.Range("E178").Value = "Gasto"
$ws.Range("F178").Value = "Tarjeta Santander"
$ws.Range("G178").Value = "Extra"
$ws.Range("K178").Value = 4856.18
$ws.Range("L178").Formula = "=L177-B178"
$ws.Range("M178").Value = 5
$ws.Range("N178").Formula = "=SUM(K178:M178)"
$ws.Range("O178").Formula = "=N178-4000"
$ws.Range("P178").Formula = "=O178-Ahorros!`$E`$4"

# --- Row 179 -------------------------------------------------------------
$ws.Range("A179").Value = 43601
$ws.Range("B179").Value = 24
$ws.Range("C179").Value = "Café Olé"
$ws.Range("D179").Value = "Golosinas"
$ws.Range("E179").Value = "Gasto"
$ws.Range("F179").Value = "Tarjeta Santander"
$ws.Range("G179").Value = "Extra"
$ws.Range("K179").Value = 4856.18
$ws.Range("L179").Formula = "=L178-B179"
$ws.Range("M179").Value = 5
$ws.Range("N179").Formula = "=SUM(K179:M179)"
$ws.Range("O179").Formula = "=N179-4000"
$ws.Range("P179").Formula = "=O179-Ahorros!`$E`$4"

# --- Row 180 -------------------------------------------------------------
$ws.Range("A180").Value = 43601
$ws.Range("B180").Value = 18
$ws.Range("C180").Value = "Rufles"
$ws.Range("D180").Value = "Golosinas"
$ws.Range("E180").Value = "Gasto"
$ws.Range("F180").Value = "Tarjeta Santander"
$ws.Range("G180").Value = "Extra"
$ws.Range("K180").Value = 4856.18
$ws.Range("L180").Formula = "=L179-B180"
$ws.Range("M180").Value = 5
$ws.Range("N180").Formula = "=SUM(K180:M180)"
$ws.Range("O180").Formula = "=N180-4000"
$ws.Range("P180").Formula = "=O180-Ahorros!`$E`$4"

# --- Row 181 -------------------------------------------------------------
$ws.Range("A181").Value = 43601
$ws.Range("B181").Value = 9
$ws.Range("C181").Value = "Dr. Pepper"
$ws.Range("D181").Value = "Golosinas"
$ws.Range("E181").Value = "Gasto"
$ws.Range("F181").Value = "Tarjeta Santander"
$ws.Range("G181").Value = "Extra"
$ws.Range("K181").Value = 4856.18
$ws.Range("L181").Formula = "=L180-B181"
$ws.Range("M181").Value = 5
$ws.Range("N181").Formula = "=SUM(K181:M181)"
$ws.Range("O181").Formula = "=N181-4000"
$ws.Range("P181").Formula = "=O181-Ahorros!`$E`$4"

# --- View state: keep the frozen header pane, move the selection to the
#     new last cell of the table (mirrors the author's final position). --
$ws.Range("Q181").Select() | Out-Null

Write-Output "edit applied"
